$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the Excel table to include the new "20-may" attendance column
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:E31"))

# Name the new column through the table header range so the table part
# (xl/tables/table1.xml) picks up the real name instead of "Column5"
$hdr = $tbl.HeaderRowRange
$hdr.Cells.Item(1, 5).Value2 = "20-may"

# Match the date-header number format/style used by the other date columns
$ws.Range("E1").NumberFormat = $ws.Range("D1").NumberFormat

# Mark attendance ("x") for column E on every row that already has marks
# in the prior two columns (rows 3 and 5 are students without marks).
$rowsWithX = 2,4,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31
foreach ($r in $rowsWithX) {
    $ws.Range("E$r").Value2 = "x"
}

# Update the sheet view to match the scrolled/selected state in the diff
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("E6").Select()
